# Auto-generated script applying numeric updates per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 40009500
$ws.Cells.Item(33, 9).Value = 55557510
$ws.Cells.Item(33, 10).Value = 28899.715
$ws.Cells.Item(33, 11).Value = 55557510
$ws.Cells.Item(33, 12).Value = 28899.715
$ws.Cells.Item(33, 13).Value = -55557281
$ws.Cells.Item(33, 14).Value = -29357.715
$ws.Cells.Item(107, 8).Value = 300.18182
$ws.Cells.Item(107, 9).Value = 190.28572
$ws.Cells.Item(107, 11).Value = 190.28572
$ws.Cells.Item(107, 13).Value = 1729.71428
$ws.Cells.Item(137, 8).Value = 1416.5588
$ws.Cells.Item(137, 9).Value = 1160.6897
$ws.Cells.Item(137, 10).Value = 2900.6
$ws.Cells.Item(137, 11).Value = 3482.0691
$ws.Cells.Item(137, 12).Value = 8701.799999999999
$ws.Cells.Item(137, 13).Value = -932.0690999999997
$ws.Cells.Item(137, 14).Value = -13801.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 663.9149
$ws.Cells.Item(2, 9).Value = 601
$ws.Cells.Item(2, 10).Value = 828.46155
$ws.Cells.Item(2, 11).Value = 601
$ws.Cells.Item(2, 12).Value = 828.46155
$ws.Cells.Item(2, 13).Value = -488
$ws.Cells.Item(2, 14).Value = -1054.46155
$ws.Cells.Item(23, 8).Value = 13636.363
$ws.Cells.Item(61, 8).Value = 1336.6154
$ws.Cells.Item(61, 9).Value = 1120.6666
$ws.Cells.Item(61, 10).Value = 1631.091
$ws.Cells.Item(61, 11).Value = 1120.6666
$ws.Cells.Item(61, 12).Value = 1631.091
$ws.Cells.Item(61, 13).Value = -908.6666
$ws.Cells.Item(61, 14).Value = -2055.091
$ws.Cells.Item(116, 8).Value = 663.9149
$ws.Cells.Item(116, 9).Value = 601
$ws.Cells.Item(116, 10).Value = 828.46155
$ws.Cells.Item(116, 11).Value = 601
$ws.Cells.Item(116, 12).Value = 828.46155
$ws.Cells.Item(116, 13).Value = 1693
$ws.Cells.Item(116, 14).Value = -5416.46155
$ws.Cells.Item(122, 8).Value = 2281.074
$ws.Cells.Item(122, 9).Value = 2208
$ws.Cells.Item(122, 11).Value = 6624
$ws.Cells.Item(122, 13).Value = -4174
$ws.Cells.Item(132, 8).Value = 1401.3684
$ws.Cells.Item(132, 9).Value = 914.6539
$ws.Cells.Item(132, 10).Value = 2455.9167
$ws.Cells.Item(132, 11).Value = 2743.9617
$ws.Cells.Item(132, 12).Value = 7367.750100000001
$ws.Cells.Item(132, 13).Value = -213.9616999999998
$ws.Cells.Item(132, 14).Value = -12427.7501
$ws.Cells.Item(136, 8).Value = 1336.6154
$ws.Cells.Item(136, 9).Value = 1120.6666
$ws.Cells.Item(136, 10).Value = 1631.091
$ws.Cells.Item(136, 11).Value = 3361.9998
$ws.Cells.Item(136, 12).Value = 4893.272999999999
$ws.Cells.Item(136, 13).Value = -811.9998000000001
$ws.Cells.Item(136, 14).Value = -9993.272999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 663.9149
$ws.Cells.Item(3, 9).Value = 601
$ws.Cells.Item(3, 10).Value = 828.46155
$ws.Cells.Item(3, 11).Value = 601
$ws.Cells.Item(3, 12).Value = 828.46155
$ws.Cells.Item(3, 13).Value = -487
$ws.Cells.Item(3, 14).Value = -1056.46155

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(19, 8).Value = 309.375
$ws.Cells.Item(19, 9).Value = 309.375
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 309.375
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = -139.375
$ws.Cells.Item(24, 8).Value = 309.375
$ws.Cells.Item(24, 9).Value = 309.375
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 309.375
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 13).Value = -139.375
$ws.Cells.Item(99, 8).Value = 2213.25
$ws.Cells.Item(99, 9).Value = 1861.3846
$ws.Cells.Item(99, 10).Value = 3738
$ws.Cells.Item(99, 11).Value = 1861.3846
$ws.Cells.Item(99, 12).Value = 3738
$ws.Cells.Item(99, 13).Value = -363.3846000000001
$ws.Cells.Item(99, 14).Value = -6734
$ws.Cells.Item(107, 8).Value = 1221.9
$ws.Cells.Item(107, 9).Value = 631.1875
$ws.Cells.Item(107, 10).Value = 3584.75
$ws.Cells.Item(107, 11).Value = 631.1875
$ws.Cells.Item(107, 12).Value = 3584.75
$ws.Cells.Item(107, 13).Value = 1288.8125
$ws.Cells.Item(107, 14).Value = -7424.75
$ws.Cells.Item(126, 8).Value = 2213.25
$ws.Cells.Item(126, 9).Value = 1861.3846
$ws.Cells.Item(126, 10).Value = 3738
$ws.Cells.Item(126, 11).Value = 5584.1538
$ws.Cells.Item(126, 12).Value = 11214
$ws.Cells.Item(126, 13).Value = -3114.1538
$ws.Cells.Item(126, 14).Value = -16154
$ws.Cells.Item(19, 14).ClearContents()
$ws.Cells.Item(24, 14).ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 852.36365
$ws.Cells.Item(5, 9).Value = 853.8570999999999
$ws.Cells.Item(5, 10).Value = 849.75
$ws.Cells.Item(5, 11).Value = 2561.5713
$ws.Cells.Item(5, 12).Value = 2549.25
$ws.Cells.Item(5, 13).Value = -2449.5713
$ws.Cells.Item(5, 14).Value = -2773.25
$ws.Cells.Item(12, 8).Value = 56.444443
$ws.Cells.Item(12, 10).Value = 62.25
$ws.Cells.Item(12, 12).Value = 186.75
$ws.Cells.Item(12, 14).Value = -532.75
$ws.Cells.Item(86, 8).Value = 448.85715
$ws.Cells.Item(86, 9).Value = 321
$ws.Cells.Item(86, 10).Value = 500
$ws.Cells.Item(86, 11).Value = 963
$ws.Cells.Item(86, 12).Value = 1500
$ws.Cells.Item(86, 13).Value = 223
$ws.Cells.Item(86, 14).Value = -3872
$ws.Cells.Item(88, 8).Value = 3199.9
$ws.Cells.Item(88, 10).Value = 3199.9
$ws.Cells.Item(88, 12).Value = 9599.700000000001
$ws.Cells.Item(88, 14).Value = -10455.7
$ws.Cells.Item(89, 8).Value = 448.85715
$ws.Cells.Item(89, 9).Value = 321
$ws.Cells.Item(89, 10).Value = 500
$ws.Cells.Item(89, 11).Value = 2889
$ws.Cells.Item(89, 12).Value = 4500
$ws.Cells.Item(89, 13).Value = 3039
$ws.Cells.Item(89, 14).Value = -16356
$ws.Cells.Item(91, 8).Value = 3199.9
$ws.Cells.Item(91, 10).Value = 3199.9
$ws.Cells.Item(91, 12).Value = 9599.700000000001
$ws.Cells.Item(91, 14).Value = -12563.7
$ws.Cells.Item(135, 8).Value = 852.36365
$ws.Cells.Item(135, 9).Value = 853.8570999999999
$ws.Cells.Item(135, 10).Value = 849.75
$ws.Cells.Item(135, 11).Value = 7684.7139
$ws.Cells.Item(135, 12).Value = 7647.75
$ws.Cells.Item(135, 13).Value = -5149.7139
$ws.Cells.Item(135, 14).Value = -12717.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 39020.89
$ws.Cells.Item(102, 9).Value = 51297.5
$ws.Cells.Item(102, 10).Value = 3944.8572
$ws.Cells.Item(102, 11).Value = 51297.5
$ws.Cells.Item(102, 12).Value = 3944.8572
$ws.Cells.Item(102, 13).Value = -49675.5
$ws.Cells.Item(102, 14).Value = -7188.8572
$ws.Cells.Item(122, 8).Value = 2694.4
$ws.Cells.Item(122, 9).Value = 2874.2727
$ws.Cells.Item(122, 10).Value = 2199.75
$ws.Cells.Item(122, 11).Value = 8622.8181
$ws.Cells.Item(122, 12).Value = 6599.25
$ws.Cells.Item(122, 13).Value = -6172.8181
$ws.Cells.Item(122, 14).Value = -11499.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(17, 8).Value = 1081.6
$ws.Cells.Item(17, 9).Value = 1002.6667
$ws.Cells.Item(17, 10).Value = 1200
$ws.Cells.Item(17, 11).Value = 1002.6667
$ws.Cells.Item(17, 12).Value = 1200
$ws.Cells.Item(17, 13).Value = -832.6667
$ws.Cells.Item(17, 14).Value = -1540
$ws.Cells.Item(22, 8).Value = 1146.1538
$ws.Cells.Item(22, 9).Value = 400
$ws.Cells.Item(22, 10).Value = 1208.3334
$ws.Cells.Item(22, 11).Value = 400
$ws.Cells.Item(22, 12).Value = 1208.3334
$ws.Cells.Item(22, 13).Value = -105
$ws.Cells.Item(22, 14).Value = -1798.3334
$ws.Cells.Item(27, 8).Value = 1146.1538
$ws.Cells.Item(27, 9).Value = 400
$ws.Cells.Item(27, 10).Value = 1208.3334
$ws.Cells.Item(27, 11).Value = 400
$ws.Cells.Item(27, 12).Value = 1208.3334
$ws.Cells.Item(27, 13).Value = -293
$ws.Cells.Item(27, 14).Value = -1422.3334
$ws.Cells.Item(40, 8).Value = 2295.3809
$ws.Cells.Item(40, 9).Value = 2066.8333
$ws.Cells.Item(40, 10).Value = 3666.6667
$ws.Cells.Item(40, 11).Value = 2066.8333
$ws.Cells.Item(40, 12).Value = 3666.6667
$ws.Cells.Item(40, 13).Value = -1930.8333
$ws.Cells.Item(40, 14).Value = -3938.6667
$ws.Cells.Item(46, 8).Value = 1375
$ws.Cells.Item(46, 9).Value = 1150
$ws.Cells.Item(46, 10).Value = 1600
$ws.Cells.Item(46, 11).Value = 1150
$ws.Cells.Item(46, 12).Value = 1600
$ws.Cells.Item(46, 13).Value = -962
$ws.Cells.Item(46, 14).Value = -1976
$ws.Cells.Item(122, 8).Value = 6054.108
$ws.Cells.Item(122, 9).Value = 6300.067
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 18900.201
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -16450.201
$ws.Cells.Item(122, 14).Value = -19900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 642.5714
$ws.Cells.Item(107, 9).Value = 277
$ws.Cells.Item(107, 10).Value = 1300.6
$ws.Cells.Item(107, 11).Value = 831
$ws.Cells.Item(107, 12).Value = 3901.8
$ws.Cells.Item(107, 13).Value = 1089
$ws.Cells.Item(107, 14).Value = -7741.799999999999
$ws.Cells.Item(132, 8).Value = 1630.0312
$ws.Cells.Item(132, 9).Value = 1260.75
$ws.Cells.Item(132, 10).Value = 1851.6
$ws.Cells.Item(132, 11).Value = 3782.25
$ws.Cells.Item(132, 12).Value = 5554.799999999999
$ws.Cells.Item(132, 13).Value = -1252.25
$ws.Cells.Item(132, 14).Value = -10614.8
$ws.Cells.Item(136, 8).Value = 5038
$ws.Cells.Item(136, 9).Value = 1104.7
$ws.Cells.Item(136, 10).Value = 14871.25
$ws.Cells.Item(136, 11).Value = 3314.1
$ws.Cells.Item(136, 12).Value = 44613.75
$ws.Cells.Item(136, 13).Value = -764.1000000000004
$ws.Cells.Item(136, 14).Value = -49713.75

Write-Host "Applied 215 value updates and 2 clears"
